$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set in_service (column E) to TRUE for rows 10-15
$ws.Range("E10:E15").Value = $true
